$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (R) to the right of the existing "2020" column (Q),
# mirroring Q's formatting for the header (row 4) and data (row 5) cells.

# Header cell: year label, same style as Q4
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

# Data cell: new data point, same style as Q5
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 42.9

# Move the selection to the new rightmost column (row 9), matching the
# original workbook's selection which tracked the last data column.
$ws.Range("R9").Select()
